# Eliminar duplicado de Carlos Aso
# Había dos registros: 'Carlos Aso' y 'Carlos Aso Miranda' con mensaje idéntico.
# Mantener 'Carlos Aso Miranda' (fila 8) y eliminar 'Carlos Aso' (fila 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find and delete the row where Nombre="Carlos" and Apellido="Aso" (exact match,
# not "Aso Miranda"), shifting the remaining rows up.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = $lastRow; $r -ge 2; $r--) {
    $nombre = $ws.Cells.Item($r, 1).Value2
    $apellido = $ws.Cells.Item($r, 2).Value2
    if ($nombre -eq "Carlos" -and $apellido -eq "Aso") {
        $ws.Rows.Item($r).Delete()
        break
    }
}
